# Generate Report for Handback
# Adds a new handback record (549a2072-f07f-4ff8-a568-550b32a85606.md) as a
# new row to each of the three tables/worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileBase = "549a2072-f07f-4ff8-a568-550b32a85606"
$mdName = "$fileBase.md"
$mdDisplayOverview = "e2e\$fileBase.md"
$commitHash = "318505e89b6bb35cdb7bf770fbd1f1bf96a33823"

$hyperlinkColor = 15570276  # BGR packed value of RGB(0x64,0x95,0xED) -> FF6495ED
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Set-HyperlinkStyle($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# Writes a plain-text value to a cell while avoiding Excel's automatic
# coercion of look-alike values (True/False/dates/numbers) into non-text
# cell types. A leading apostrophe forces text entry; re-applying the
# "Normal" style afterwards clears the resulting quote-prefix indicator so
# the cell format matches a plain, unstyled text cell.
function Set-TextValue($range, [string]$text) {
    $range.Value2 = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$null = $loOverview.ListRows.Add()

$wsOverview.Range("A3").Value2 = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$mdName", $null, $null, $mdDisplayOverview)
Set-HyperlinkStyle $wsOverview.Range("B3")
$wsOverview.Range("C3").Value2 = ".md"
$wsOverview.Range("E3").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value2 = "2016-10-27 05:58:01"
$wsOverview.Range("G3").NumberFormat = $dateFormat

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$null = $loZh.ListRows.Add()

$zhXlf = "$fileBase.$commitHash.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$mdName", $null, $null, $mdName)
Set-HyperlinkStyle $wsZh.Range("A3")
$wsZh.Range("B3").Value2 = ".md"
$wsZh.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value2 = "e2e"
$wsZh.Range("E3").Value2 = "ht"
Set-TextValue $wsZh.Range("F3") "True"
$wsZh.Range("G3").Value2 = $zhXlf
Set-TextValue $wsZh.Range("H3") "2016-10-27 05:57:47"
$wsZh.Range("H3").NumberFormat = $dateFormat
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$commitHash/e2e/$mdName", $null, $null, $mdName)
Set-HyperlinkStyle $wsZh.Range("I3")
$wsZh.Range("J3").Value2 = $zhXlf
Set-TextValue $wsZh.Range("K3") "2016-10-27 05:58:38"
$wsZh.Range("K3").NumberFormat = $dateFormat
Set-TextValue $wsZh.Range("L3") ""
Set-TextValue $wsZh.Range("M3") "True"
Set-TextValue $wsZh.Range("N3") ""
Set-TextValue $wsZh.Range("O3") "False"
Set-TextValue $wsZh.Range("P3") ""

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$null = $loDe.ListRows.Add()

$deXlf = "$fileBase.$commitHash.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$mdName", $null, $null, $mdName)
Set-HyperlinkStyle $wsDe.Range("A3")
$wsDe.Range("B3").Value2 = ".md"
$wsDe.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value2 = "e2e"
$wsDe.Range("E3").Value2 = "ht"
$wsDe.Range("F3").Value2 = "'True"
$wsDe.Range("G3").Value2 = $deXlf
$wsDe.Range("H3").Value2 = "2016-10-27 05:58:01"
$wsDe.Range("H3").NumberFormat = $dateFormat
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$commitHash/e2e/$mdName", $null, $null, $mdName)
Set-HyperlinkStyle $wsDe.Range("I3")
$wsDe.Range("J3").Value2 = $deXlf
$wsDe.Range("K3").Value2 = "2016-10-27 05:58:55"
$wsDe.Range("K3").NumberFormat = $dateFormat
$wsDe.Range("L3").Value2 = "'"
$wsDe.Range("M3").Value2 = "'True"
$wsDe.Range("N3").Value2 = "'"
$wsDe.Range("O3").Value2 = "'False"
$wsDe.Range("P3").Value2 = "'"
